$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = '27/12/2025 00:49'
$ws.Cells.Item(4, 3).Value = 542
$ws.Cells.Item(4, 4).Value = 'Conhecimentos Específicos'
$ws.Cells.Item(4, 5).Value = 'Gestão da Manutenção e Confiabilidade'
$ws.Cells.Item(4, 6).Value = 0

$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = '27/12/2025 00:58'
$ws.Cells.Item(5, 3).Value = 1105
$ws.Cells.Item(5, 4).Value = 'Estatística'
$ws.Cells.Item(5, 5).Value = 'Variáveis Aleatórias e Distribuições Discretas'
$ws.Cells.Item(5, 6).Value = 1

$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = '27/12/2025 02:59'
$ws.Cells.Item(6, 3).Value = '980'
$ws.Cells.Item(6, 4).Value = 'Inglês'
$ws.Cells.Item(6, 5).Value = 'Interpretação de Texto'
$ws.Cells.Item(6, 6).Value = 0

$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = '27/12/2025 03:10'
$ws.Cells.Item(7, 3).Value = '980'
$ws.Cells.Item(7, 4).Value = 'Inglês'
$ws.Cells.Item(7, 5).Value = 'Interpretação de Texto'
$ws.Cells.Item(7, 6).Value = 0

$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = '27/12/2025 03:11'
$ws.Cells.Item(8, 3).Value = 979
$ws.Cells.Item(8, 4).Value = 'Inglês'
$ws.Cells.Item(8, 5).Value = 'Interpretação de Texto'
$ws.Cells.Item(8, 6).Value = 0

$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = '27/12/2025 03:11'
$ws.Cells.Item(9, 3).Value = 893
$ws.Cells.Item(9, 4).Value = 'Inglês'
$ws.Cells.Item(9, 5).Value = 'Interpretação de Texto'
$ws.Cells.Item(9, 6).Value = 0

$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = '27/12/2025 03:11'
$ws.Cells.Item(10, 3).Value = 892
$ws.Cells.Item(10, 4).Value = 'Inglês'
$ws.Cells.Item(10, 5).Value = 'Semantic'
$ws.Cells.Item(10, 6).Value = 0

$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = '27/12/2025 03:11'
$ws.Cells.Item(11, 3).Value = 891
$ws.Cells.Item(11, 4).Value = 'Inglês'
$ws.Cells.Item(11, 5).Value = 'Interpretação de Texto'
$ws.Cells.Item(11, 6).Value = 1

$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = '27/12/2025 03:21'
$ws.Cells.Item(12, 3).Value = 1208
$ws.Cells.Item(12, 4).Value = 'Estatística'
$ws.Cells.Item(12, 5).Value = 'Regressão Linear Múltipla'
$ws.Cells.Item(12, 6).Value = 0

$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = '27/12/2025 03:28'
$ws.Cells.Item(13, 3).Value = 630
$ws.Cells.Item(13, 4).Value = 'Conhecimentos Específicos'
$ws.Cells.Item(13, 5).Value = 'Gestão da Qualidade'
$ws.Cells.Item(13, 6).Value = 1

$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = '27/12/2025 03:31'
$ws.Cells.Item(14, 3).Value = 1203
$ws.Cells.Item(14, 4).Value = 'Estatística'
$ws.Cells.Item(14, 5).Value = 'Regressão Linear Múltipla'
$ws.Cells.Item(14, 6).Value = 0

$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = '27/12/2025 03:40'
$ws.Cells.Item(15, 3).Value = 451
$ws.Cells.Item(15, 4).Value = 'Conhecimentos Específicos'
$ws.Cells.Item(15, 5).Value = 'Gestão de Estoques'
$ws.Cells.Item(15, 6).Value = 1
